$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume data (and the Cosmos / InjectiveProtocol
# row swap) as scraped by the GitHub Actions workflow.
$updates = @(
    @{ Cell = "D2"; Value = "42.580.87"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -0.86%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "2.528.83"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -0.20%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  +0.00%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "317.53"; ForceText = $true },
    @{ Cell = "D6"; Value = "94.77"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -6.16%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  -0.94%  "; ForceText = $false },
    @{ Cell = "E8"; Value = "  -0.13%  "; ForceText = $false },
    @{ Cell = "E9"; Value = "  -3.57%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "35.98"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -4.12%  "; ForceText = $false },
    @{ Cell = "E11"; Value = "  -1.40%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "0.113"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -0.19%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "7.54"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -2.68%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "2.914.17"; ForceText = $false },
    @{ Cell = "E14"; Value = "  -0.27%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "2.527.81"; ForceText = $false },
    @{ Cell = "E15"; Value = "  +2.10%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "15.38"; ForceText = $true },
    @{ Cell = "D17"; Value = "0.846"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -3.03%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "42.554.93"; ForceText = $false },
    @{ Cell = "E18"; Value = "  -0.86%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "12.99"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -1.90%  "; ForceText = $false },
    @{ Cell = "E20"; Value = "  +0.88%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "0.0₃0960"; ForceText = $false },
    @{ Cell = "E21"; Value = "  -2.85%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "69.99"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -2.30%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "250.98"; ForceText = $true },
    @{ Cell = "E23"; Value = "  -1.31%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "2.96"; ForceText = $true },
    @{ Cell = "E24"; Value = "  +0.75%  "; ForceText = $false },
    @{ Cell = "E25"; Value = "  -2.26%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "26.46"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -3.06%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "0.996"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -0.61%  "; ForceText = $false },
    @{ Cell = "E28"; Value = "  +2.87%  "; ForceText = $false },
    @{ Cell = "B29"; Value = "InjectiveProtocol"; ForceText = $false },
    @{ Cell = "C29"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; ForceText = $false },
    @{ Cell = "D29"; Value = "39.11"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +0.99%  "; ForceText = $false },
    @{ Cell = "B30"; Value = "Cosmos"; ForceText = $false },
    @{ Cell = "C30"; Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; ForceText = $false },
    @{ Cell = "D30"; Value = "10.16"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -2.07%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "6.04"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -1.86%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "154.85"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -1.70%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "19.22"; ForceText = $true },
    @{ Cell = "E33"; Value = "  +4.68%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "2.11"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -0.32%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "3.27"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -0.73%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "0.0784"; ForceText = $true },
    @{ Cell = "E36"; Value = "  -1.60%  "; ForceText = $false },
    @{ Cell = "E37"; Value = "  -0.74%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "0.111"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -3.80%  "; ForceText = $false },
    @{ Cell = "E39"; Value = "  -0.75%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "23.65"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -0.51%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "2.31"; ForceText = $true },
    @{ Cell = "E41"; Value = "  +10.02%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  +0.37%  "; ForceText = $false },
    @{ Cell = "E43"; Value = "  -2.46%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "0.0300"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -1.62%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "3.27"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -6.23%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "2.006.01"; ForceText = $false },
    @{ Cell = "E46"; Value = "  -1.91%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "85.15"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -1.30%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "8.78"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -2.41%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "2.767.77"; ForceText = $false },
    @{ Cell = "E49"; Value = "  -0.45%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "74.00"; ForceText = $true },
    @{ Cell = "E50"; Value = "  +1.07%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "102.42"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -0.73%  "; ForceText = $false }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    if ($update.ForceText) {
        # Force text storage so numeric-looking strings (e.g. "317.53")
        # are not auto-converted to actual numbers by Excel.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $update.Value
}
